$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Vega Central Mapocho de Santiago - Acelga"
# block (rows 991-992), pushing all existing data down by two rows. This mirrors the
# weekly refresh: a new week (2023-12-07) of price data is added and the oldest rows
# fall off the bottom of the A1:R range, growing the sheet from R1049 to R1051.
$ws.Rows("991:992").Insert()

# New row 991: Acelga, "Primera" grade, week of 2023-12-07.
$ws.Cells.Item(991, 1).Value = 9
$ws.Cells.Item(991, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(991, 3).Value = "Metropolitana"
$ws.Cells.Item(991, 4).Value = [DateTime]"2023-12-07"
$ws.Cells.Item(991, 5).Value = 13
$ws.Cells.Item(991, 6).Value = 100112009
$ws.Cells.Item(991, 7).Value = "Acelga"
$ws.Cells.Item(991, 8).Value = "Sin especificar"
$ws.Cells.Item(991, 9).Value = "Primera"
$ws.Cells.Item(991, 10).Value = 160
$ws.Cells.Item(991, 11).Value = 22000
$ws.Cells.Item(991, 12).Value = 24000
$ws.Cells.Item(991, 13).Value = 23000
$ws.Cells.Item(991, 14).Value = "$/docena de atados"
$ws.Cells.Item(991, 15).Value = "Región Metropolitana"
$ws.Cells.Item(991, 16).Value = 7667
$ws.Cells.Item(991, 17).Value = 3
$ws.Cells.Item(991, 18).Value = "Hortaliza"

# New row 992: Acelga, "Segunda" grade, same week.
$ws.Cells.Item(992, 1).Value = 9
$ws.Cells.Item(992, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(992, 3).Value = "Metropolitana"
$ws.Cells.Item(992, 4).Value = [DateTime]"2023-12-07"
$ws.Cells.Item(992, 5).Value = 13
$ws.Cells.Item(992, 6).Value = 100112009
$ws.Cells.Item(992, 7).Value = "Acelga"
$ws.Cells.Item(992, 8).Value = "Sin especificar"
$ws.Cells.Item(992, 9).Value = "Segunda"
$ws.Cells.Item(992, 10).Value = 70
$ws.Cells.Item(992, 11).Value = 18000
$ws.Cells.Item(992, 12).Value = 18000
$ws.Cells.Item(992, 13).Value = 18000
$ws.Cells.Item(992, 14).Value = "$/docena de atados"
$ws.Cells.Item(992, 15).Value = "Región Metropolitana"
$ws.Cells.Item(992, 16).Value = 6000
$ws.Cells.Item(992, 17).Value = 3
$ws.Cells.Item(992, 18).Value = "Hortaliza"
